$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, border, centered) from H1 to I1:J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-23
$data = @{
    2  = @(7, 8)
    3  = @(1, 3)
    4  = @(1, 6)
    5  = @(1, 6)
    6  = @(4, 6)
    7  = @(7, 8)
    8  = @(6, 9)
    9  = @(8, 9)
    10 = @(6, 9)
    11 = @(1, 5)
    12 = @(8, 9)
    13 = @(1, 6)
    14 = @(1, 4)
    15 = @(1, 6)
    16 = @(1, 3)
    17 = @(1, 5)
    18 = @(1, 6)
    19 = @(1, 6)
    20 = @(1, 4)
    21 = @(1, 4)
    22 = @(1, 2)
    23 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}

Write-Output "Added columns I0 and IF"
